# Auto-generated edit script: reorders D, M, N, O, P, Q, R, S, T columns across rows 2-26
# per a row permutation (weekly data re-association), per commit 'Fruta / hortaliza, semanal'.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (values sourced from original row 16)
$ws.Cells.Item(2, 4).Value = [datetime]"2021-02-18"
$ws.Cells.Item(2, 13).Value = 50
$ws.Cells.Item(2, 14).Value = 15000
$ws.Cells.Item(2, 15).Value = 15000
$ws.Cells.Item(2, 16).Value = 15000
$ws.Cells.Item(2, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(2, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(2, 19).Value = 1000
$ws.Cells.Item(2, 20).Value = 15

# Row 3 (values sourced from original row 21)
$ws.Cells.Item(3, 4).Value = [datetime]"2021-05-04"
$ws.Cells.Item(3, 13).Value = 45
$ws.Cells.Item(3, 14).Value = 14000
$ws.Cells.Item(3, 15).Value = 14000
$ws.Cells.Item(3, 16).Value = 14000
$ws.Cells.Item(3, 17).Value = "$/caja 14 kilos granel"
$ws.Cells.Item(3, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(3, 19).Value = 1000
$ws.Cells.Item(3, 20).Value = 14

# Row 4 (values sourced from original row 19)
$ws.Cells.Item(4, 4).Value = [datetime]"2022-03-10"
$ws.Cells.Item(4, 13).Value = 75
$ws.Cells.Item(4, 14).Value = 15000
$ws.Cells.Item(4, 15).Value = 15000
$ws.Cells.Item(4, 16).Value = 15000
$ws.Cells.Item(4, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(4, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(4, 19).Value = 1071
$ws.Cells.Item(4, 20).Value = 14

# Row 5 (values sourced from original row 10)
$ws.Cells.Item(5, 4).Value = [datetime]"2022-01-31"
$ws.Cells.Item(5, 13).Value = 54
$ws.Cells.Item(5, 14).Value = 20000
$ws.Cells.Item(5, 15).Value = 20000
$ws.Cells.Item(5, 16).Value = 20000
$ws.Cells.Item(5, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(5, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(5, 19).Value = 1333
$ws.Cells.Item(5, 20).Value = 15

# Row 6 (values sourced from original row 26)
$ws.Cells.Item(6, 4).Value = [datetime]"2021-03-16"
$ws.Cells.Item(6, 13).Value = 50
$ws.Cells.Item(6, 14).Value = 12000
$ws.Cells.Item(6, 15).Value = 12000
$ws.Cells.Item(6, 16).Value = 12000
$ws.Cells.Item(6, 17).Value = "$/caja 14 kilos granel"
$ws.Cells.Item(6, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(6, 19).Value = 857
$ws.Cells.Item(6, 20).Value = 14

# Row 7 (values sourced from original row 15)
$ws.Cells.Item(7, 4).Value = [datetime]"2021-04-27"
$ws.Cells.Item(7, 13).Value = 36
$ws.Cells.Item(7, 14).Value = 14000
$ws.Cells.Item(7, 15).Value = 14000
$ws.Cells.Item(7, 16).Value = 14000
$ws.Cells.Item(7, 17).Value = "$/caja 14 kilos granel"
$ws.Cells.Item(7, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(7, 19).Value = 1000
$ws.Cells.Item(7, 20).Value = 14

# Row 8 (values sourced from original row 14)
$ws.Cells.Item(8, 4).Value = [datetime]"2022-02-22"
$ws.Cells.Item(8, 13).Value = 54
$ws.Cells.Item(8, 14).Value = 14000
$ws.Cells.Item(8, 15).Value = 14000
$ws.Cells.Item(8, 16).Value = 14000
$ws.Cells.Item(8, 17).Value = "$/caja 14 kilos granel"
$ws.Cells.Item(8, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(8, 19).Value = 1000
$ws.Cells.Item(8, 20).Value = 14

# Row 9 (values sourced from original row 13)
$ws.Cells.Item(9, 4).Value = [datetime]"2021-04-26"
$ws.Cells.Item(9, 13).Value = 68
$ws.Cells.Item(9, 14).Value = 14000
$ws.Cells.Item(9, 15).Value = 14000
$ws.Cells.Item(9, 16).Value = 14000
$ws.Cells.Item(9, 17).Value = "$/caja 14 kilos granel"
$ws.Cells.Item(9, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(9, 19).Value = 1000
$ws.Cells.Item(9, 20).Value = 14

# Row 10 (values sourced from original row 22)
$ws.Cells.Item(10, 4).Value = [datetime]"2021-03-23"
$ws.Cells.Item(10, 13).Value = 45
$ws.Cells.Item(10, 14).Value = 13000
$ws.Cells.Item(10, 15).Value = 13000
$ws.Cells.Item(10, 16).Value = 13000
$ws.Cells.Item(10, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(10, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(10, 19).Value = 929
$ws.Cells.Item(10, 20).Value = 14

# Row 11 (values sourced from original row 25)
$ws.Cells.Item(11, 4).Value = [datetime]"2021-03-05"
$ws.Cells.Item(11, 13).Value = 56
$ws.Cells.Item(11, 14).Value = 13000
$ws.Cells.Item(11, 15).Value = 13000
$ws.Cells.Item(11, 16).Value = 13000
$ws.Cells.Item(11, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(11, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(11, 19).Value = 929
$ws.Cells.Item(11, 20).Value = 14

# Row 12 (values sourced from original row 24)
$ws.Cells.Item(12, 4).Value = [datetime]"2021-05-06"
$ws.Cells.Item(12, 13).Value = 50
$ws.Cells.Item(12, 14).Value = 14000
$ws.Cells.Item(12, 15).Value = 14000
$ws.Cells.Item(12, 16).Value = 14000
$ws.Cells.Item(12, 17).Value = "$/caja 14 kilos granel"
$ws.Cells.Item(12, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(12, 19).Value = 1000
$ws.Cells.Item(12, 20).Value = 14

# Row 13 (values sourced from original row 3)
$ws.Cells.Item(13, 4).Value = [datetime]"2021-03-04"
$ws.Cells.Item(13, 13).Value = 80
$ws.Cells.Item(13, 14).Value = 12000
$ws.Cells.Item(13, 15).Value = 12000
$ws.Cells.Item(13, 16).Value = 12000
$ws.Cells.Item(13, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(13, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(13, 19).Value = 800
$ws.Cells.Item(13, 20).Value = 15

# Row 14 (values sourced from original row 6)
$ws.Cells.Item(14, 4).Value = [datetime]"2021-03-15"
$ws.Cells.Item(14, 13).Value = 85
$ws.Cells.Item(14, 14).Value = 12000
$ws.Cells.Item(14, 15).Value = 12000
$ws.Cells.Item(14, 16).Value = 12000
$ws.Cells.Item(14, 17).Value = "$/caja 14 kilos granel"
$ws.Cells.Item(14, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(14, 19).Value = 857
$ws.Cells.Item(14, 20).Value = 14

# Row 15 (values sourced from original row 2)
$ws.Cells.Item(15, 4).Value = [datetime]"2021-02-25"
$ws.Cells.Item(15, 13).Value = 60
$ws.Cells.Item(15, 14).Value = 14000
$ws.Cells.Item(15, 15).Value = 14000
$ws.Cells.Item(15, 16).Value = 14000
$ws.Cells.Item(15, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(15, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(15, 19).Value = 1000
$ws.Cells.Item(15, 20).Value = 14

# Row 16 (values sourced from original row 7)
$ws.Cells.Item(16, 4).Value = [datetime]"2022-01-27"
$ws.Cells.Item(16, 13).Value = 85
$ws.Cells.Item(16, 14).Value = 19000
$ws.Cells.Item(16, 15).Value = 20000
$ws.Cells.Item(16, 16).Value = 19529
$ws.Cells.Item(16, 17).Value = "$/caja 14 kilos granel"
$ws.Cells.Item(16, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(16, 19).Value = 1395
$ws.Cells.Item(16, 20).Value = 14

# Row 17 (values sourced from original row 5)
$ws.Cells.Item(17, 4).Value = [datetime]"2022-02-24"
$ws.Cells.Item(17, 13).Value = 70
$ws.Cells.Item(17, 14).Value = 14000
$ws.Cells.Item(17, 15).Value = 14000
$ws.Cells.Item(17, 16).Value = 14000
$ws.Cells.Item(17, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(17, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(17, 19).Value = 1000
$ws.Cells.Item(17, 20).Value = 14

# Row 18 (values sourced from original row 9)
$ws.Cells.Item(18, 4).Value = [datetime]"2022-01-24"
$ws.Cells.Item(18, 13).Value = 50
$ws.Cells.Item(18, 14).Value = 22500
$ws.Cells.Item(18, 15).Value = 22500
$ws.Cells.Item(18, 16).Value = 22500
$ws.Cells.Item(18, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(18, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(18, 19).Value = 1500
$ws.Cells.Item(18, 20).Value = 15

# Row 19 (values sourced from original row 4)
$ws.Cells.Item(19, 4).Value = [datetime]"2021-02-11"
$ws.Cells.Item(19, 13).Value = 60
$ws.Cells.Item(19, 14).Value = 15000
$ws.Cells.Item(19, 15).Value = 15000
$ws.Cells.Item(19, 16).Value = 15000
$ws.Cells.Item(19, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(19, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(19, 19).Value = 1000
$ws.Cells.Item(19, 20).Value = 15

# Row 20 (values sourced from original row 18)
$ws.Cells.Item(20, 4).Value = [datetime]"2021-04-28"
$ws.Cells.Item(20, 13).Value = 56
$ws.Cells.Item(20, 14).Value = 14000
$ws.Cells.Item(20, 15).Value = 14000
$ws.Cells.Item(20, 16).Value = 14000
$ws.Cells.Item(20, 17).Value = "$/caja 14 kilos granel"
$ws.Cells.Item(20, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(20, 19).Value = 1000
$ws.Cells.Item(20, 20).Value = 14

# Row 21 (values sourced from original row 12)
$ws.Cells.Item(21, 4).Value = [datetime]"2022-03-07"
$ws.Cells.Item(21, 13).Value = 56
$ws.Cells.Item(21, 14).Value = 17000
$ws.Cells.Item(21, 15).Value = 17000
$ws.Cells.Item(21, 16).Value = 17000
$ws.Cells.Item(21, 17).Value = "$/caja 14 kilos empedrada"
$ws.Cells.Item(21, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(21, 19).Value = 1214
$ws.Cells.Item(21, 20).Value = 14

# Row 22 (values sourced from original row 17)
$ws.Cells.Item(22, 4).Value = [datetime]"2021-02-12"
$ws.Cells.Item(22, 13).Value = 70
$ws.Cells.Item(22, 14).Value = 15000
$ws.Cells.Item(22, 15).Value = 15000
$ws.Cells.Item(22, 16).Value = 15000
$ws.Cells.Item(22, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(22, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(22, 19).Value = 1000
$ws.Cells.Item(22, 20).Value = 15

# Row 23 (values sourced from original row 11)
$ws.Cells.Item(23, 4).Value = [datetime]"2021-02-15"
$ws.Cells.Item(23, 13).Value = 45
$ws.Cells.Item(23, 14).Value = 12000
$ws.Cells.Item(23, 15).Value = 12000
$ws.Cells.Item(23, 16).Value = 12000
$ws.Cells.Item(23, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(23, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(23, 19).Value = 800
$ws.Cells.Item(23, 20).Value = 15

# Row 24 (values sourced from original row 23)
$ws.Cells.Item(24, 4).Value = [datetime]"2021-04-30"
$ws.Cells.Item(24, 13).Value = 48
$ws.Cells.Item(24, 14).Value = 14000
$ws.Cells.Item(24, 15).Value = 14000
$ws.Cells.Item(24, 16).Value = 14000
$ws.Cells.Item(24, 17).Value = "$/caja 14 kilos granel"
$ws.Cells.Item(24, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(24, 19).Value = 1000
$ws.Cells.Item(24, 20).Value = 14

# Row 25 (values sourced from original row 8)
$ws.Cells.Item(25, 4).Value = [datetime]"2021-05-07"
$ws.Cells.Item(25, 13).Value = 60
$ws.Cells.Item(25, 14).Value = 14000
$ws.Cells.Item(25, 15).Value = 14000
$ws.Cells.Item(25, 16).Value = 14000
$ws.Cells.Item(25, 17).Value = "$/caja 14 kilos granel"
$ws.Cells.Item(25, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(25, 19).Value = 1000
$ws.Cells.Item(25, 20).Value = 14

# Row 26 (values sourced from original row 20)
$ws.Cells.Item(26, 4).Value = [datetime]"2021-04-29"
$ws.Cells.Item(26, 13).Value = 65
$ws.Cells.Item(26, 14).Value = 14000
$ws.Cells.Item(26, 15).Value = 14000
$ws.Cells.Item(26, 16).Value = 14000
$ws.Cells.Item(26, 17).Value = "$/caja 14 kilos granel"
$ws.Cells.Item(26, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(26, 19).Value = 1000
$ws.Cells.Item(26, 20).Value = 14

